$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Refresh the "datetimeFigureOut" date placeholders (slide master + every
#    slide layout) from 11/09/2018 -> 27/09/2018.
# ---------------------------------------------------------------------------
function Update-DateShapes($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -eq "11/09/2018") {
                $sh.TextFrame.TextRange.Text = "27/09/2018"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-DateShapes $master.CustomLayouts.Item($li)
}

# ---------------------------------------------------------------------------
# 2. Slide 1: widen/move the "Ogre" rounded-rectangle callout and retitle it
#    to "OgrDB".
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)

# a:off x="6224631" -> x="5343787" ; a:ext cx="4320165" -> cx="5201009"
# (y / cy are unchanged). Left/Width are expressed in points (1pt = 12700 EMU).
$shape.Left = 420.7706604003906
$shape.Width = 409.5282897949219

$shape.TextFrame.TextRange.Text = "OgrDB"
